$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.392.16"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "1.805.43"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'327.49"
$ws.Range("E5").Value = "  -3.08%  "
$ws.Range("D6").Value = "'0.9993"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.4453"
$ws.Range("E7").Value = "  +6.10%  "
$ws.Range("D8").Value = "'0.3737"
$ws.Range("E8").Value = "  +6.41%  "
$ws.Range("D9").Value = "'44.74"
$ws.Range("E9").Value = "  -2.36%  "
$ws.Range("D10").Value = "'1.149"
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("D11").Value = "'0.07514"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "'22.57"
$ws.Range("E12").Value = "  -1.24%  "
$ws.Range("D13").Value = "'1.001"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "'7.702"
$ws.Range("E14").Value = "  +5.43%  "
$ws.Range("D15").Value = "'6.297"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "1.801.19"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "'0.00001094"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "'0.06789"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("D19").Value = "'80.87"
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").Value = "'0.9994"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "'17.47"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").Value = "'6.325"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").Value = "28.352.78"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").Value = "'11.80"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("D25").Value = "'2.409"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").Value = "'20.48"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").Value = "'153.42"
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("D28").Value = "'2.351"
$ws.Range("E28").Value = "  -5.06%  "
$ws.Range("D29").Value = "2.003.67"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").Value = "'132.51"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "'1.253"
$ws.Range("E31").Value = "  -4.92%  "
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("D33").Value = "'5.823"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("D34").Value = "'0.09334"
$ws.Range("E34").Value = "  +2.75%  "
$ws.Range("D35").Value = "'0.2287"
$ws.Range("E35").Value = "  +5.92%  "
$ws.Range("D36").Value = "'12.12"
$ws.Range("E36").Value = "  -2.31%  "
$ws.Range("D37").Value = "'0.06341"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "'0.02321"
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "'5.166"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6574"
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("D42").Value = "'1.461"
$ws.Range("E42").Value = "  -3.50%  "
$ws.Range("D43").Value = "'8.180"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "'0.9990"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'14.08"
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("D46").Value = "'0.6069"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").Value = "'3.796"
$ws.Range("E47").Value = "  -2.10%  "
$ws.Range("D48").Value = "'128.30"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").Value = "'2.033"
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("D50").Value = "'0.07107"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").Value = "'1.157"
$ws.Range("E51").Value = "  -2.26%  "
